$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value into a cell while keeping it stored as text
# (the source data keeps every column, including numeric-looking ones,
# as text cells).
function Set-TextValue {
    param($addr, $value)
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value = $value
}

# 1) Straightforward nombre_aides / montant_total revisions for existing
#    rows (counts/amounts revised upward with the 2020-08-27 data refresh).
$updates = @{
    2  = @{ C = "810";  D = "1827687.79" }
    4  = @{ C = "1030"; D = "3668453.47" }
    6  = @{ C = "667";  D = "2161907.78" }
    7  = @{ C = "17";   D = "37498.41" }
    30 = @{ C = "568";  D = "2353516.89" }
    45 = @{ C = "396";  D = "1065299.43" }
    48 = @{ C = "435";  D = "1534361.40" }
    51 = @{ C = "3819"; D = "8806141.47" }
    56 = @{ C = "4194"; D = "13329755.27" }
    61 = @{ C = "25";   D = "91177.00" }
    74 = @{ C = "387";  D = "979509.70" }
    76 = @{ C = "935";  D = "3275240.26" }
    79 = @{ C = "36";   D = "156180.27" }
    86 = @{ C = "223";  D = "526169.00" }
    88 = @{ C = "512";  D = "1822869.71" }
    89 = @{ C = "188";  D = "556380.81" }
    90 = @{ C = "17";   D = "57000.00" }
    91 = @{ C = "8";    D = "24670.00" }
    92 = @{ C = "668";  D = "1630584.94" }
}

foreach ($row in $updates.Keys) {
    foreach ($col in $updates[$row].Keys) {
        Set-TextValue "$col$row" $updates[$row][$col]
    }
}

# 2) A new "Société créée de fait" (code 22) breakdown row for the
#    Provence-Alpes-Côte d'Azur region is added, inserted just above the
#    former row 93 ("Société en nom collectif"), pushing it and every
#    category-juridique row below it down by one.
$ws.Rows.Item(93).Insert()

Set-TextValue "A93" "Fonds de solidarité"
Set-TextValue "B93" "VOLET2"
Set-TextValue "C93" "3"
Set-TextValue "D93" "7571.16"
Set-TextValue "E93" "93"
Set-TextValue "F93" "Provence-Alpes-Côte d'Azur"
Set-TextValue "G93" "22"
Set-TextValue "H93" "Société créée de fait"

# 3) The rows that used to be 93-98 (now 94-99) keep their region/category
#    pairing but several of their nombre_aides / montant_total totals were
#    also revised as part of the same data refresh.
$shiftedUpdates = @{
    95 = @{ C = "1113"; D = "3782549.80" }
    97 = @{ C = "1044"; D = "3243097.76" }
    99 = @{ C = "52";   D = "203835.52" }
}

foreach ($row in $shiftedUpdates.Keys) {
    foreach ($col in $shiftedUpdates[$row].Keys) {
        Set-TextValue "$col$row" $shiftedUpdates[$row][$col]
    }
}
